# Update cryptos list (Price/Volume(1h) columns) with latest scraped values.
# For D-column (Price) values that look like plain decimal numbers, force the
# cell's number format to Text ("@") first so Excel keeps the exact original
# string (e.g. "71.40") instead of silently coercing it into a float (71.4)
# and dropping the trailing zero / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.652.12'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '3.250.43'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.95'
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.23'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  +1.31%  '
$ws.Range("D9").Value = '3.249.15'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("E11").Value = '  -2.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.412'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '3.799.58'
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.68'
$ws.Range("E15").Value = '  -3.06%  '
$ws.Range("D16").Value = '67.648.45'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '3.212.87'
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.77'
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.53'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '394.70'
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.58'
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.40'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.517'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.61'
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.96'
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.56'
$ws.Range("E31").Value = '  -4.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.67'
$ws.Range("E32").Value = '  -1.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.02'
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("E34").Value = '  -2.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.16'
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  -4.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.89'
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.61'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.809'
$ws.Range("E40").Value = '  -3.61%  '
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.48'
$ws.Range("E42").Value = '  -4.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.48'
$ws.Range("E43").Value = '  -7.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0688'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.69'
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("D46").Value = '2.616.48'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.74'
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '334.91'
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0279'
$ws.Range("E49").Value = '  -2.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.33'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E51").Value = '  -0.29%  '
